# Generate Report for Handoff
# Update "Latest Handoff Datetime" (column D, row 5 -> file
# 10598507-7115-4d94-9882-232ac59c87d2.md, status "Ready for handoff")
# on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D5").Value = "2016-03-09 16:31:33"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D5").Value = "2016-03-09 16:31:43"
